# Updated symbol list on Wed Dec 14 06:11:09 UTC 2022 with GitHub Actions
# Refreshes the "Price" (column D) and "Hora" (column G) columns for each
# coin row on the active sheet. Values are written with a leading
# apostrophe so Excel keeps them as text (matching the original inlineStr
# cell type) instead of auto-converting numeric-looking text to a Number;
# the Style reset afterwards clears the transient "quote prefix" format
# flag that the apostrophe entry leaves behind, so no unintended
# formatting change is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'274.87"
$ws.Range("G2").Value = "'6"
$ws.Range("D2").Style = "Normal"
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").Value = "'23.02"
$ws.Range("G3").Value = "'6"
$ws.Range("D3").Style = "Normal"
$ws.Range("G3").Style = "Normal"
$ws.Range("G4").Value = "'6"
$ws.Range("G4").Style = "Normal"
$ws.Range("D5").Value = "'0.06268"
$ws.Range("G5").Value = "'6"
$ws.Range("D5").Style = "Normal"
$ws.Range("G5").Style = "Normal"
$ws.Range("D6").Value = "'3.664"
$ws.Range("G6").Value = "'6"
$ws.Range("D6").Style = "Normal"
$ws.Range("G6").Style = "Normal"
$ws.Range("D7").Value = "'6.666"
$ws.Range("G7").Value = "'6"
$ws.Range("D7").Style = "Normal"
$ws.Range("G7").Style = "Normal"
$ws.Range("D8").Value = "'1.349"
$ws.Range("G8").Value = "'6"
$ws.Range("D8").Style = "Normal"
$ws.Range("G8").Style = "Normal"
$ws.Range("D9").Value = "'0.8320"
$ws.Range("G9").Value = "'6"
$ws.Range("D9").Style = "Normal"
$ws.Range("G9").Style = "Normal"
$ws.Range("D10").Value = "'0.01375"
$ws.Range("G10").Value = "'6"
$ws.Range("D10").Style = "Normal"
$ws.Range("G10").Style = "Normal"
$ws.Range("D11").Value = "'0.1635"
$ws.Range("G11").Value = "'6"
$ws.Range("D11").Style = "Normal"
$ws.Range("G11").Style = "Normal"
$ws.Range("D12").Value = "'0.08332"
$ws.Range("G12").Value = "'6"
$ws.Range("D12").Style = "Normal"
$ws.Range("G12").Style = "Normal"
$ws.Range("D13").Value = "'0.03417"
$ws.Range("G13").Value = "'6"
$ws.Range("D13").Style = "Normal"
$ws.Range("G13").Style = "Normal"
$ws.Range("D14").Value = "'0.03086"
$ws.Range("G14").Value = "'6"
$ws.Range("D14").Style = "Normal"
$ws.Range("G14").Style = "Normal"
$ws.Range("D15").Value = "'0.09307"
$ws.Range("G15").Value = "'6"
$ws.Range("D15").Style = "Normal"
$ws.Range("G15").Style = "Normal"
$ws.Range("D16").Value = "'3.839"
$ws.Range("G16").Value = "'6"
$ws.Range("D16").Style = "Normal"
$ws.Range("G16").Style = "Normal"
$ws.Range("D17").Value = "'0.001633"
$ws.Range("G17").Value = "'6"
$ws.Range("D17").Style = "Normal"
$ws.Range("G17").Style = "Normal"
$ws.Range("D18").Value = "'0.04768"
$ws.Range("G18").Value = "'6"
$ws.Range("D18").Style = "Normal"
$ws.Range("G18").Style = "Normal"
$ws.Range("D19").Value = "'0.006391"
$ws.Range("G19").Value = "'6"
$ws.Range("D19").Style = "Normal"
$ws.Range("G19").Style = "Normal"
$ws.Range("D20").Value = "'0.005682"
$ws.Range("G20").Value = "'6"
$ws.Range("D20").Style = "Normal"
$ws.Range("G20").Style = "Normal"
$ws.Range("D21").Value = "'0.001089"
$ws.Range("G21").Value = "'6"
$ws.Range("D21").Style = "Normal"
$ws.Range("G21").Style = "Normal"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("G22").Value = "'6"
$ws.Range("D22").Style = "Normal"
$ws.Range("G22").Style = "Normal"
$ws.Range("D23").Value = "'3.716"
$ws.Range("G23").Value = "'6"
$ws.Range("D23").Style = "Normal"
$ws.Range("G23").Style = "Normal"
$ws.Range("G24").Value = "'6"
$ws.Range("G24").Style = "Normal"
$ws.Range("D25").Value = "'0.3346"
$ws.Range("G25").Value = "'6"
$ws.Range("D25").Style = "Normal"
$ws.Range("G25").Style = "Normal"
$ws.Range("G26").Value = "'6"
$ws.Range("G26").Style = "Normal"
$ws.Range("D27").Value = "'0.0002679"
$ws.Range("G27").Value = "'6"
$ws.Range("D27").Style = "Normal"
$ws.Range("G27").Style = "Normal"
$ws.Range("G28").Value = "'6"
$ws.Range("G28").Style = "Normal"
$ws.Range("G29").Value = "'6"
$ws.Range("G29").Style = "Normal"
$ws.Range("G30").Value = "'6"
$ws.Range("G30").Style = "Normal"
$ws.Range("G31").Value = "'6"
$ws.Range("G31").Style = "Normal"
$ws.Range("G32").Value = "'6"
$ws.Range("G32").Style = "Normal"
$ws.Range("G33").Value = "'6"
$ws.Range("G33").Style = "Normal"
$ws.Range("G34").Value = "'6"
$ws.Range("G34").Style = "Normal"
$ws.Range("G35").Value = "'6"
$ws.Range("G35").Style = "Normal"
$ws.Range("G36").Value = "'6"
$ws.Range("G36").Style = "Normal"
$ws.Range("G37").Value = "'6"
$ws.Range("G37").Style = "Normal"
$ws.Range("G38").Value = "'6"
$ws.Range("G38").Style = "Normal"
$ws.Range("G39").Value = "'6"
$ws.Range("G39").Style = "Normal"
$ws.Range("D40").Value = "'0.04702"
$ws.Range("G40").Value = "'6"
$ws.Range("D40").Style = "Normal"
$ws.Range("G40").Style = "Normal"
$ws.Range("D41").Value = "'0.007047"
$ws.Range("G41").Value = "'6"
$ws.Range("D41").Style = "Normal"
$ws.Range("G41").Style = "Normal"
$ws.Range("D42").Value = "'0.1163"
$ws.Range("G42").Value = "'6"
$ws.Range("D42").Style = "Normal"
$ws.Range("G42").Style = "Normal"
$ws.Range("D43").Value = "'0.003519"
$ws.Range("G43").Value = "'6"
$ws.Range("D43").Style = "Normal"
$ws.Range("G43").Style = "Normal"
$ws.Range("D44").Value = "'0.01217"
$ws.Range("G44").Value = "'6"
$ws.Range("D44").Style = "Normal"
$ws.Range("G44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006253"
$ws.Range("G45").Value = "'6"
$ws.Range("D45").Style = "Normal"
$ws.Range("G45").Style = "Normal"
$ws.Range("G46").Value = "'6"
$ws.Range("G46").Style = "Normal"
$ws.Range("G47").Value = "'6"
$ws.Range("G47").Style = "Normal"
$ws.Range("D48").Value = "'0.7963"
$ws.Range("G48").Value = "'6"
$ws.Range("D48").Style = "Normal"
$ws.Range("G48").Style = "Normal"
$ws.Range("D49").Value = "'0.03617"
$ws.Range("G49").Value = "'6"
$ws.Range("D49").Style = "Normal"
$ws.Range("G49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002299"
$ws.Range("G50").Value = "'6"
$ws.Range("D50").Style = "Normal"
$ws.Range("G50").Style = "Normal"
$ws.Range("D51").Value = "'0.01240"
$ws.Range("G51").Value = "'6"
$ws.Range("D51").Style = "Normal"
$ws.Range("G51").Style = "Normal"
